$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 10; $r++) {
    $ws.Range("B${r}:G${r}").Value = "N/A"
}

# The "% of nodes required to take over network" column (E) previously
# carried a Percent number format (style index 1); now that every cell in
# it holds the literal text "N/A", drop that formatting back to Normal so
# the unused Percent style/font can fall away on save.
$ws.Range("E2:E10").Style = "Normal"

$ws.Range("J8").Select()
